# Updates the "cryptos" price/volume table to the latest scraped values.
# Row 34/35 additionally swap Fetch.AI <-> Aptos (their rank order changed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.938.48'
$ws.Range("E2").Value = '  -2.54%  '

$ws.Range("D3").Value = '3.474.83'
$ws.Range("E3").Value = '  +0.82%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '''582.24'
$ws.Range("E5").Value = '  -1.51%  '

$ws.Range("D6").Value = '''173.04'
$ws.Range("E6").Value = '  -3.14%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = '''0.600'
$ws.Range("E8").Value = '  -1.68%  '

$ws.Range("D9").Value = '3.475.53'
$ws.Range("E9").Value = '  +0.96%  '

$ws.Range("E10").Value = '  -6.03%  '

$ws.Range("D11").Value = '''6.86'
$ws.Range("E11").Value = '  -1.57%  '

$ws.Range("D12").Value = '''0.411'
$ws.Range("E12").Value = '  -3.76%  '

$ws.Range("D13").Value = '4.074.46'
$ws.Range("E13").Value = '  +0.79%  '

$ws.Range("E14").Value = '  +0.56%  '

$ws.Range("D15").Value = '''29.93'
$ws.Range("E15").Value = '  -6.39%  '

$ws.Range("D16").Value = '66.064.95'
$ws.Range("E16").Value = '  -2.29%  '

$ws.Range("E17").Value = '  -3.16%  '

$ws.Range("D18").Value = '3.469.15'
$ws.Range("E18").Value = '  +0.74%  '

$ws.Range("D19").Value = '''5.92'
$ws.Range("E19").Value = '  -3.37%  '

$ws.Range("D20").Value = '''13.88'
$ws.Range("E20").Value = '  -0.78%  '

$ws.Range("D21").Value = '''366.66'
$ws.Range("E21").Value = '  -5.65%  '

$ws.Range("E22").Value = '  -1.54%  '

$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("D24").Value = '''72.17'
$ws.Range("E24").Value = '  +1.10%  '

$ws.Range("D25").Value = '''0.535'
$ws.Range("E25").Value = '  +0.48%  '

$ws.Range("E26").Value = '  +4.44%  '

$ws.Range("D27").Value = '''9.57'
$ws.Range("E27").Value = '  -6.41%  '

$ws.Range("E28").Value = '  +2.30%  '

$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.12%  '

$ws.Range("D30").Value = '''23.90'
$ws.Range("E30").Value = '  +2.87%  '

$ws.Range("D31").Value = '''5.76'
$ws.Range("E31").Value = '  -4.77%  '

$ws.Range("E32").Value = '  -2.98%  '

$ws.Range("D33").Value = '''1.00'
$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = '''7.11'
$ws.Range("E34").Value = '  -1.39%  '

$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").Value = '''1.29'
$ws.Range("E35").Value = '  -6.77%  '

$ws.Range("D36").Value = '''1.53'
$ws.Range("E36").Value = '  -1.52%  '

$ws.Range("D37").Value = '''159.40'
$ws.Range("E37").Value = '  -0.92%  '

$ws.Range("D38").Value = '''29.28'
$ws.Range("E38").Value = '  +13.23%  '

$ws.Range("D39").Value = '''0.889'
$ws.Range("E39").Value = '  +0.89%  '

$ws.Range("D40").Value = '2.809.42'
$ws.Range("E40").Value = '  +4.22%  '

$ws.Range("E41").Value = '  -5.42%  '

$ws.Range("E42").Value = '  -6.83%  '

$ws.Range("D43").Value = '''6.45'
$ws.Range("E43").Value = '  -2.81%  '

$ws.Range("E44").Value = '  -3.54%  '

$ws.Range("D45").Value = '''0.0682'
$ws.Range("E45").Value = '  -4.35%  '

$ws.Range("D46").Value = '''40.12'
$ws.Range("E46").Value = '  -2.49%  '

$ws.Range("E47").Value = '  -7.12%  '

$ws.Range("D48").Value = '''0.0288'
$ws.Range("E48").Value = '  -2.63%  '

$ws.Range("D49").Value = '''313.63'
$ws.Range("E49").Value = '  -2.95%  '

$ws.Range("E50").Value = '  -2.59%  '

$ws.Range("E51").Value = '  -2.61%  '
